# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets carry identical data, so the same per-row updates apply to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 131
    3  = 45
    6  = 124
    7  = 1218
    8  = 1504
    10 = 374
    12 = 138
    13 = 164
    17 = 292
    19 = 1705
    22 = 172
    23 = 651
    26 = 4105
    29 = 255
    33 = 458
    35 = 211
    36 = 48
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
